$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = -7.069
$ws.Range("D12").Value = -7.203
$ws.Range("D15").Value = -8.209
$ws.Range("D27").Value = -7.792
$ws.Range("D28").Value = -8.147000000000002
$ws.Range("D31").Value = -7.770000000000001
$ws.Range("D32").Value = -7.052000000000001
$ws.Range("D36").Value = -7.657000000000001
$ws.Range("D38").Value = -7.706
$ws.Range("D46").Value = -7.859000000000002
$ws.Range("D54").Value = -7.878
$ws.Range("D55").Value = -8.028
$ws.Range("D56").Value = -8.321000000000002
$ws.Range("D67").Value = -7.556999999999999
$ws.Range("D69").Value = -7.274000000000001
$ws.Range("D72").Value = -7.398000000000001
$ws.Range("D73").Value = -8.104000000000001
$ws.Range("D83").Value = -7.849000000000001
$ws.Range("D86").Value = -8.245999999999999
$ws.Range("D91").Value = -7.636
$ws.Range("D93").Value = -7.203
$ws.Range("D99").Value = -8.074999999999999
$ws.Range("D104").Value = -7.587999999999999
$ws.Range("D105").Value = -7.956
